$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 120945.87
$ws.Range("I138").Value = 6074.3335
$ws.Range("J138").Value = 134032.5
$ws.Range("K138").Value = 18223.0005
$ws.Range("L138").Value = 402097.5
$ws.Range("M138").Value = -13083.0005
$ws.Range("N138").Value = -412377.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 295
$ws.Range("I2").Value = 295
$ws.Range("K2").Value = 295
$ws.Range("M2").Value = -182
$ws.Range("H32").Value = 24700.158
$ws.Range("I32").Value = 20715.424
$ws.Range("K32").Value = 20715.424
$ws.Range("M32").Value = -20428.424
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H61").Value = 2958216.8
$ws.Range("I61").Value = 5001985
$ws.Range("K61").Value = 5001985
$ws.Range("M61").Value = -5001773
$ws.Range("H74").Value = 4084013.2
$ws.Range("I74").Value = 4763522
$ws.Range("K74").Value = 4763522
$ws.Range("M74").Value = -4762648
$ws.Range("H77").Value = 4084013.2
$ws.Range("I77").Value = 4763522
$ws.Range("K77").Value = 23817610
$ws.Range("M77").Value = -23813242
$ws.Range("H116").Value = 295
$ws.Range("I116").Value = 295
$ws.Range("K116").Value = 295
$ws.Range("M116").Value = 1999
$ws.Range("H136").Value = 2958216.8
$ws.Range("I136").Value = 5001985
$ws.Range("K136").Value = 15005955
$ws.Range("M136").Value = -15003405
$ws.Range("H137").Value = 95181.82000000001
$ws.Range("J137").Value = 94700
$ws.Range("L137").Value = 94700
$ws.Range("N137").Value = -104900

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 295
$ws.Range("I3").Value = 295
$ws.Range("K3").Value = 295
$ws.Range("M3").Value = -181
$ws.Range("H99").Value = 3389
$ws.Range("I99").Value = 2625.0557
$ws.Range("J99").Value = 6139.2
$ws.Range("K99").Value = 2625.0557
$ws.Range("L99").Value = 6139.2
$ws.Range("M99").Value = -1127.0557
$ws.Range("N99").Value = -9135.200000000001
$ws.Range("H134").Value = 3288.7144
$ws.Range("I134").Value = 2960.318
$ws.Range("K134").Value = 8880.954000000002
$ws.Range("M134").Value = -6345.954000000002
$ws.Range("H140").Value = 165251.45
$ws.Range("J140").Value = 165251.45
$ws.Range("L140").Value = 165251.45
$ws.Range("N140").Value = -175611.45

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4052.1785
$ws.Range("I31").Value = 3798.647
$ws.Range("K31").Value = 3798.647
$ws.Range("M31").Value = -3503.647
$ws.Range("H34").Value = 4052.1785
$ws.Range("I34").Value = 3798.647
$ws.Range("K34").Value = 3798.647
$ws.Range("M34").Value = -3596.647
$ws.Range("H52").Value = 98999.39999999999
$ws.Range("I52").Value = 50000
$ws.Range("J52").Value = 111249.25
$ws.Range("K52").Value = 50000
$ws.Range("L52").Value = 111249.25
$ws.Range("M52").Value = -49706
$ws.Range("N52").Value = -111837.25
$ws.Range("H58").Value = 5093.067
$ws.Range("I58").Value = 6559.2
$ws.Range("K58").Value = 6559.2
$ws.Range("M58").Value = -6356.2
$ws.Range("H132").Value = 5027.7666
$ws.Range("I132").Value = 4149.2104
$ws.Range("K132").Value = 12447.6312
$ws.Range("M132").Value = -9917.6312
$ws.Range("H136").Value = 5093.067
$ws.Range("I136").Value = 6559.2
$ws.Range("K136").Value = 19677.6
$ws.Range("M136").Value = -17127.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 839.6
$ws.Range("I8").Value = 839.6
$ws.Range("K8").Value = 2518.8
$ws.Range("M8").Value = -2379.8
$ws.Range("H122").Value = 2599.4
$ws.Range("J122").Value = 3332.3333
$ws.Range("L122").Value = 29990.9997
$ws.Range("N122").Value = -34890.9997
$ws.Range("H131").Value = 9999999
$ws.Range("I131").Value = 9999999
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 29999997
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -29994957
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 5628.75
$ws.Range("I136").Value = 5628.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 16886.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -11786.25
$ws.Range("N136").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5653.5
$ws.Range("I113").Value = 4784.2
$ws.Range("K113").Value = 4784.2
$ws.Range("M113").Value = -2614.2
$ws.Range("H132").Value = 5115.694
$ws.Range("I132").Value = 5253.8335
$ws.Range("K132").Value = 15761.5005
$ws.Range("M132").Value = -13231.5005
$ws.Range("H135").Value = 106000
$ws.Range("J135").Value = 106000
$ws.Range("L135").Value = 106000
$ws.Range("N135").Value = -116140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 10241.571
$ws.Range("I93").Value = 11968.2
$ws.Range("J93").Value = 5925
$ws.Range("K93").Value = 11968.2
$ws.Range("L93").Value = 5925
$ws.Range("M93").Value = -10720.2
$ws.Range("N93").Value = -8421
$ws.Range("H100").Value = 7288.3125
$ws.Range("J100").Value = 4800
$ws.Range("L100").Value = 4800
$ws.Range("N100").Value = -5882
$ws.Range("H132").Value = 4195.7427
$ws.Range("I132").Value = 4221.9585
$ws.Range("K132").Value = 12665.8755
$ws.Range("M132").Value = -10135.8755
$ws.Range("H136").Value = 6071.1787
$ws.Range("I136").Value = 5217.087
$ws.Range("K136").Value = 15651.261
$ws.Range("M136").Value = -13101.261
$ws.Range("H137").Value = 83500
$ws.Range("J137").Value = 99666.664
$ws.Range("L137").Value = 99666.664
$ws.Range("N137").Value = -109866.664

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 14999.333
$ws.Range("I5").Value = 14998
$ws.Range("J5").Value = 15000
$ws.Range("K5").Value = 14998
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = -14886
$ws.Range("N5").Value = -15224
$ws.Range("H126").Value = 2420.25
$ws.Range("I126").Value = 2406.4736
$ws.Range("J126").Value = 2472.6
$ws.Range("K126").Value = 7219.4208
$ws.Range("L126").Value = 7417.799999999999
$ws.Range("M126").Value = -4749.4208
$ws.Range("N126").Value = -12357.8
$ws.Range("H131").Value = 94274.664
$ws.Range("J131").Value = 92999.60000000001
$ws.Range("L131").Value = 92999.60000000001
$ws.Range("N131").Value = -103079.6
$ws.Range("H136").Value = 15854.75
$ws.Range("I136").Value = 24520.666
$ws.Range("K136").Value = 73561.99800000001
$ws.Range("M136").Value = -71011.99800000001
$ws.Range("H137").Value = 86918.8
$ws.Range("J137").Value = 86918.8
$ws.Range("L137").Value = 86918.8
$ws.Range("N137").Value = -97118.8
$ws.Range("H139").Value = 249579
$ws.Range("I139").Value = 849995
$ws.Range("K139").Value = 849995
$ws.Range("M139").Value = -844855
